$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix polyAIsolationProtocol column (G) value for every data row: was
# "NEBNextPoly(A)E7490L", should be "E7420L".
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 7).Value = "E7420L"
}

# Widen column G (polyAIsolationProtocol) now that formatting needed fixing.
$ws.Columns.Item(7).ColumnWidth = 27.5

# A few stray blank rows exist further down the sheet.
$ws.Rows.Item(39).RowHeight = 15
$ws.Rows.Item(40).RowHeight = 15
$ws.Rows.Item(43).RowHeight = 15

# Restore the view position/selection as last left by the editor.
$ws.Application.Goto($ws.Range("A8"), $false)
$ws.Range("G25").Select()
